$d = $word.ActiveDocument

# 1) Table formatting: normalize the table's preferred width (2500.0% -> 2500,
#    i.e. still 50%) and mark the first row as a repeating header row.
$t = $d.Tables.Item(1)
$t.PreferredWidthType = 2   # wdPreferredWidthPercent
$t.PreferredWidth = 125     # 125 * 20 = 2500 fiftieths-of-a-percent = 50%
$t.Rows.Item(1).HeadingFormat = $true

# 2) Trim the "Sketch the t-distribution using the t-distribution applet."
#    sentence from the P-value instruction bullet.
$d.Content.Find.Execute(
    "Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Find the P-value and compare it to the level of significance.",
    2)
